$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated log write mode: refreshed simulated run metrics (run_time, max_er, iter 0-19)

$ws.Range("C2").Value = 0.2159323692321777
$ws.Range("E2").Value = 116.7257070795258
$ws.Range("F2").Value = 0.003403769951593227
$ws.Range("G2").Value = 0.003403769951593227
$ws.Range("H2").Value = 0.002717440152790865
$ws.Range("I2").Value = 0.002717440152790865
$ws.Range("J2").Value = 0.002717440152790865
$ws.Range("K2").Value = 0.002717440152790865
$ws.Range("L2").Value = 0.002652220947650199
$ws.Range("M2").Value = 0.002416901744621881
$ws.Range("N2").Value = 0.002394072392101993
$ws.Range("O2").Value = 0.002394072392101993
$ws.Range("P2").Value = 0.002394072392101993
$ws.Range("Q2").Value = 0.002394072392101993
$ws.Range("R2").Value = 0.002367892145530188
$ws.Range("S2").Value = 0.002367892145530188
$ws.Range("T2").Value = 0.002326878983193323
$ws.Range("U2").Value = 0.002317315338941076
$ws.Range("V2").Value = 0.002315901694191443
$ws.Range("W2").Value = 0.0023042798363505
$ws.Range("X2").Value = 0.002294312593934464
$ws.Range("Y2").Value = 0.002275354913830912

$ws.Range("C3").Value = 0.2181625366210938
$ws.Range("E3").Value = 116.8324859888871
$ws.Range("F3").Value = 0.003400105573240949
$ws.Range("G3").Value = 0.003124351940380178
$ws.Range("H3").Value = 0.002808730397142713
$ws.Range("I3").Value = 0.002808730397142713
$ws.Range("J3").Value = 0.002808730397142713
$ws.Range("K3").Value = 0.002808730397142713
$ws.Range("L3").Value = 0.002808730397142713
$ws.Range("M3").Value = 0.002808730397142713
$ws.Range("N3").Value = 0.002563003279799241
$ws.Range("O3").Value = 0.002539884647674264
$ws.Range("P3").Value = 0.002539884647674264
$ws.Range("Q3").Value = 0.002536549760987125
$ws.Range("R3").Value = 0.002433149990425824
$ws.Range("S3").Value = 0.002391001170361593
$ws.Range("T3").Value = 0.002355483154031241
$ws.Range("U3").Value = 0.002299806719234361
$ws.Range("V3").Value = 0.002286380550403617
$ws.Range("W3").Value = 0.002286380550403617
$ws.Range("X3").Value = 0.002286380550403617
$ws.Range("Y3").Value = 0.00227743637405238

$ws.Range("C4").Value = 0.187507152557373
$ws.Range("E4").Value = 127.1094564371542
$ws.Range("F4").Value = 0.003301287609747925
$ws.Range("G4").Value = 0.002939996241491163
$ws.Range("H4").Value = 0.002831265394267857
$ws.Range("I4").Value = 0.002831265394267857
$ws.Range("J4").Value = 0.002831265394267857
$ws.Range("K4").Value = 0.002831265394267857
$ws.Range("L4").Value = 0.002831265394267857
$ws.Range("M4").Value = 0.002753111770238617
$ws.Range("N4").Value = 0.002753111770238617
$ws.Range("O4").Value = 0.002596065274188398
$ws.Range("P4").Value = 0.002596065274188398
$ws.Range("Q4").Value = 0.002596065274188398
$ws.Range("R4").Value = 0.002596065274188398
$ws.Range("S4").Value = 0.002596065274188398
$ws.Range("T4").Value = 0.002596065274188398
$ws.Range("U4").Value = 0.002552064526612757
$ws.Range("V4").Value = 0.002552064526612757
$ws.Range("W4").Value = 0.00251247124298552
$ws.Range("X4").Value = 0.002477767182010803
$ws.Range("Y4").Value = 0.002477767182010803

$ws.Range("C5").Value = 0.253443717956543
$ws.Range("E5").Value = 122.8535501446277
$ws.Range("F5").Value = 0.003522423937944115
$ws.Range("G5").Value = 0.00284824391404815
$ws.Range("H5").Value = 0.00284824391404815
$ws.Range("I5").Value = 0.00284824391404815
$ws.Range("J5").Value = 0.00283408032370466
$ws.Range("K5").Value = 0.002833473895416832
$ws.Range("L5").Value = 0.002650296359011921
$ws.Range("M5").Value = 0.002650296359011921
$ws.Range("N5").Value = 0.002650296359011921
$ws.Range("O5").Value = 0.002465197328437512
$ws.Range("P5").Value = 0.002465197328437512
$ws.Range("Q5").Value = 0.002465197328437512
$ws.Range("R5").Value = 0.002465197328437512
$ws.Range("S5").Value = 0.002465197328437512
$ws.Range("T5").Value = 0.002465197328437512
$ws.Range("U5").Value = 0.002465197328437512
$ws.Range("V5").Value = 0.002465197328437512
$ws.Range("W5").Value = 0.002407104814015487
$ws.Range("X5").Value = 0.002407104814015487
$ws.Range("Y5").Value = 0.002394806045704242

$ws.Range("C6").Value = 0.3339650630950928
$ws.Range("E6").Value = 121.7811056707287
$ws.Range("F6").Value = 0.00335144582154993
$ws.Range("G6").Value = 0.003095131316536265
$ws.Range("H6").Value = 0.003095131316536265
$ws.Range("I6").Value = 0.002994597531247123
$ws.Range("J6").Value = 0.002841638606137122
$ws.Range("K6").Value = 0.002841638606137122
$ws.Range("L6").Value = 0.002643742622574181
$ws.Range("M6").Value = 0.002643742622574181
$ws.Range("N6").Value = 0.002643742622574181
$ws.Range("O6").Value = 0.002643742622574181
$ws.Range("P6").Value = 0.002626676172517284
$ws.Range("Q6").Value = 0.002527187489995519
$ws.Range("R6").Value = 0.002527187489995519
$ws.Range("S6").Value = 0.002527187489995519
$ws.Range("T6").Value = 0.002496954951870062
$ws.Range("U6").Value = 0.00241620036380106
$ws.Range("V6").Value = 0.002403964811074245
$ws.Range("W6").Value = 0.002403964811074245
$ws.Range("X6").Value = 0.002394952530739048
$ws.Range("Y6").Value = 0.002373900695335842

$ws.Range("C7").Value = 0.3253564834594727
$ws.Range("E7").Value = 116.9860913417488
$ws.Range("F7").Value = 0.003420623059458673
$ws.Range("G7").Value = 0.002960281829183094
$ws.Range("H7").Value = 0.002960281829183094
$ws.Range("I7").Value = 0.002891106893836556
$ws.Range("J7").Value = 0.002739501626383855
$ws.Range("K7").Value = 0.002635084850179494
$ws.Range("L7").Value = 0.002635084850179494
$ws.Range("M7").Value = 0.002635084850179494
$ws.Range("N7").Value = 0.00258919448397924
$ws.Range("O7").Value = 0.00258919448397924
$ws.Range("P7").Value = 0.002458844329240025
$ws.Range("Q7").Value = 0.002458844329240025
$ws.Range("R7").Value = 0.002417267183923206
$ws.Range("S7").Value = 0.002417267183923206
$ws.Range("T7").Value = 0.002370940501961818
$ws.Range("U7").Value = 0.002370940501961818
$ws.Range("V7").Value = 0.002334372705954486
$ws.Range("W7").Value = 0.002334372705954486
$ws.Range("X7").Value = 0.002309911127221282
$ws.Range("Y7").Value = 0.002280430630443445

$ws.Range("C8").Value = 0.2296862602233887
$ws.Range("E8").Value = 118.1248554165995
$ws.Range("F8").Value = 0.003383280333885364
$ws.Range("G8").Value = 0.003050988193930496
$ws.Range("H8").Value = 0.002695292200116728
$ws.Range("I8").Value = 0.002695292200116728
$ws.Range("J8").Value = 0.002695292200116728
$ws.Range("K8").Value = 0.002695292200116728
$ws.Range("L8").Value = 0.002571740267029343
$ws.Range("M8").Value = 0.002571740267029343
$ws.Range("N8").Value = 0.002565953775809897
$ws.Range("O8").Value = 0.002565953775809897
$ws.Range("P8").Value = 0.002496029561069249
$ws.Range("Q8").Value = 0.002496029561069249
$ws.Range("R8").Value = 0.002410316585637506
$ws.Range("S8").Value = 0.002410316585637506
$ws.Range("T8").Value = 0.002395182235641097
$ws.Range("U8").Value = 0.002358396293455417
$ws.Range("V8").Value = 0.002345746883884221
$ws.Range("W8").Value = 0.002307812219402272
$ws.Range("X8").Value = 0.002307074446373617
$ws.Range("Y8").Value = 0.002302628760557495

$ws.Range("C9").Value = 0.2472198009490967
$ws.Range("E9").Value = 124.5414572782702
$ws.Range("F9").Value = 0.003521271639873229
$ws.Range("G9").Value = 0.003161100990172108
$ws.Range("H9").Value = 0.00301289634147385
$ws.Range("I9").Value = 0.00301289634147385
$ws.Range("J9").Value = 0.00301289634147385
$ws.Range("K9").Value = 0.002870590622596173
$ws.Range("L9").Value = 0.002870590622596173
$ws.Range("M9").Value = 0.002870590622596173
$ws.Range("N9").Value = 0.002870590622596173
$ws.Range("O9").Value = 0.002685979719631523
$ws.Range("P9").Value = 0.002685979719631523
$ws.Range("Q9").Value = 0.002685979719631523
$ws.Range("R9").Value = 0.002598742843779938
$ws.Range("S9").Value = 0.002575741926290457
$ws.Range("T9").Value = 0.002540515389460481
$ws.Range("U9").Value = 0.002540515389460481
$ws.Range("V9").Value = 0.002518737014139496
$ws.Range("W9").Value = 0.002457657693245784
$ws.Range("X9").Value = 0.002427708718874662
$ws.Range("Y9").Value = 0.002427708718874662

$ws.Range("C10").Value = 0.1927549839019775
$ws.Range("E10").Value = 116.8764947291183
$ws.Range("F10").Value = 0.003454616380507159
$ws.Range("G10").Value = 0.002981128491030768
$ws.Range("H10").Value = 0.002616598422188515
$ws.Range("I10").Value = 0.002616598422188515
$ws.Range("J10").Value = 0.002616598422188515
$ws.Range("K10").Value = 0.002616598422188515
$ws.Range("L10").Value = 0.002616598422188515
$ws.Range("M10").Value = 0.002616598422188515
$ws.Range("N10").Value = 0.002616598422188515
$ws.Range("O10").Value = 0.002540017864780841
$ws.Range("P10").Value = 0.002540017864780841
$ws.Range("Q10").Value = 0.002540017864780841
$ws.Range("R10").Value = 0.002506209792586036
$ws.Range("S10").Value = 0.002404686803845011
$ws.Range("T10").Value = 0.002403499027412209
$ws.Range("U10").Value = 0.002389121389009779
$ws.Range("V10").Value = 0.002361236115734019
$ws.Range("W10").Value = 0.002299954768663382
$ws.Range("X10").Value = 0.002284529951669286
$ws.Range("Y10").Value = 0.002278294244232326

$ws.Range("C11").Value = 0.2738447189331055
$ws.Range("E11").Value = 122.29647979109
$ws.Range("F11").Value = 0.003522423937944115
$ws.Range("G11").Value = 0.003171828982680271
$ws.Range("H11").Value = 0.002838111942227483
$ws.Range("I11").Value = 0.002838111942227483
$ws.Range("J11").Value = 0.002838111942227483
$ws.Range("K11").Value = 0.002757164521265021
$ws.Range("L11").Value = 0.002742421086112757
$ws.Range("M11").Value = 0.002575507828547649
$ws.Range("N11").Value = 0.002575507828547649
$ws.Range("O11").Value = 0.002575507828547649
$ws.Range("P11").Value = 0.002575507828547649
$ws.Range("Q11").Value = 0.002556501891866888
$ws.Range("R11").Value = 0.002545969992047033
$ws.Range("S11").Value = 0.002430994001781009
$ws.Range("T11").Value = 0.002430994001781009
$ws.Range("U11").Value = 0.002430994001781009
$ws.Range("V11").Value = 0.002430994001781009
$ws.Range("W11").Value = 0.002418129053648262
$ws.Range("X11").Value = 0.002392101169221136
$ws.Range("Y11").Value = 0.002383946974485184
